$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.272327238179451;  C = 0.3048912486333797;  D = 189.6080260415259; E = 13.86384647080068; F = 0; G = 207.0490909991394 }
    3 = @{ B = 0.1169995834814548; C = 0.3048912486333797;  D = 3.223369029078222;  E = 13.86384647080068; F = 0; G = 17.50910633199374 }
    4 = @{ B = 0.6545652718822623; C = 1.626987699542094;   D = 3.223369029078222;  E = 13.86384647080068; F = 1; G = 19.36876847130326 }
    5 = @{ B = 0.01253208636536152; C = 0.002658071450198252; D = 0.7210945179870265; E = 13.86384647080068; F = 1; G = 14.60013114660327 }
    6 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987; F = 0; G = 6.15379541431027 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    $ws.Range("B$row").Value = $rowData.B
    $ws.Range("C$row").Value = $rowData.C
    $ws.Range("D$row").Value = $rowData.D
    $ws.Range("E$row").Value = $rowData.E
    $ws.Range("F$row").Value = $rowData.F
    $ws.Range("G$row").Value = $rowData.G
}

$wb.Save()
